$d = $word.ActiveDocument

# The first paragraph reads "This is a Microsoft word document." and is
# a single run. The diff appends " (", "Changed main" and ")" as three
# additional runs after the existing run (so the visible text becomes
# "This is a Microsoft word document. (Changed main)").
$para = $d.Paragraphs(1).Range

# Whole-paragraph range, excluding the trailing paragraph mark.
$full = $d.Range($para.Start, $para.End)

# Rebuild the paragraph's content via InsertXML so the new text lands in
# its own separate <w:r> runs (matching the diff) instead of being
# merged into the existing run.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r><w:t>This is a Microsoft word document.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r><w:t>Changed main</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml)
